$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 21,4
$data[0,0] = 6
$data[0,1] = 4
$data[0,2] = 4
$data[0,3] = 16
$data[1,0] = 6
$data[1,1] = 12
$data[1,2] = 5
$data[1,3] = 8
$data[2,0] = 3
$data[2,1] = 13
$data[2,2] = 5
$data[2,3] = 7
$data[3,0] = 3
$data[3,1] = 14
$data[3,2] = 5
$data[3,3] = 6
$data[4,0] = 4
$data[4,1] = 8
$data[4,2] = 5
$data[4,3] = 12
$data[5,0] = 4
$data[5,1] = 16
$data[5,2] = 2
$data[5,3] = 4
$data[6,0] = 2
$data[6,1] = 8
$data[6,2] = 3
$data[6,3] = 12
$data[7,0] = 4
$data[7,1] = 6
$data[7,2] = 2
$data[7,3] = 14
$data[8,0] = 8
$data[8,1] = 13
$data[8,2] = 7
$data[8,3] = 7
$data[9,0] = 5
$data[9,1] = 7
$data[9,2] = 4
$data[9,3] = 13
$data[10,0] = 6
$data[10,1] = 7
$data[10,2] = 4
$data[10,3] = 13
$data[11,0] = 1
$data[11,1] = 5
$data[11,2] = 2
$data[11,3] = 15
$data[12,0] = 3
$data[12,1] = 13
$data[12,2] = 5
$data[12,3] = 7
$data[13,0] = 3
$data[13,1] = 0
$data[13,2] = 4
$data[13,3] = 20
$data[14,0] = 7
$data[14,1] = 8
$data[14,2] = 6
$data[14,3] = 12
$data[15,0] = 6
$data[15,1] = 5
$data[15,2] = 8
$data[15,3] = 15
$data[16,0] = 4
$data[16,1] = 17
$data[16,2] = 2
$data[16,3] = 3
$data[17,0] = 3
$data[17,1] = 8
$data[17,2] = 2
$data[17,3] = 12
$data[18,0] = 7
$data[18,1] = 14
$data[18,2] = 5
$data[18,3] = 6
$data[19,0] = 6
$data[19,1] = 12
$data[19,2] = 5
$data[19,3] = 8
$data[20,0] = 2
$data[20,1] = 16
$data[20,2] = 3
$data[20,3] = 4

$ws.Range("A1598:D1618").Value = $data

$ws.Range("A1619").Select()